# Update the CDA Logical model (PlayingEntity) "Metadata" sheet for ST.r2b:
#  - bump Version and Date values
#  - insert a new "Jurisdiction" property row (after "Contact"), pushing
#    Description/Purpose/Copyright/... down by one row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3: Property="Version")
$ws.Range("B3").Value2 = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8: Property="Date")
$ws.Range("B8").Value2 = "2025-10-29T22:15:57+01:00"

# Insert a new row after the "Contact" row (row 10) for the new
# "Jurisdiction" property.
$ws.Rows.Item(11).EntireRow.Insert()

# Match the formatting of the surrounding data rows (border/alignment),
# then fill in the new property/value pair.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value2 = "Jurisdiction"
$ws.Range("B11").Value2 = ""
